$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ALC_values = @{
    "H5" = 53.166668
    "I5" = 57.2
    "J5" = 33
    "K5" = 57.2
    "L5" = 33
    "M5" = 57.8
    "N5" = -263
    "H28" = 134388.27
    "I28" = 182756.73
    "K28" = 182756.73
    "M28" = -182271.73
    "H40" = 3256.4866
    "I40" = 1865.3334
    "J40" = 3703.6428
    "K40" = 1865.3334
    "L40" = 3703.6428
    "M40" = -1690.3334
    "N40" = -4053.6428
    "H41" = 1466.75
    "J41" = 232.33333
    "L41" = 232.33333
    "N41" = -1112.33333
    "H74" = 6181.5
    "I74" = 6377.8
    "K74" = 6377.8
    "M74" = -5441.8
    "H77" = 6181.5
    "I77" = 6377.8
    "K77" = 31889
    "M77" = -27209
    "H86" = 120376090
    "I86" = 100005920
    "J86" = 178576580
    "K86" = 100005920
    "L86" = 178576580
    "M86" = -100004797
    "N86" = -178578826
    "H89" = 120376090
    "I89" = 100005920
    "J89" = 178576580
    "K89" = 500029600
    "L89" = 892882900
    "M89" = -500023984
    "N89" = -892894132
    "H106" = 5130514
    "I106" = 5130514
    "K106" = 5130514
    "M106" = -5129883
    "H107" = 1316.4445
    "J107" = 675.25
    "L107" = 675.25
    "N107" = -4515.25
    "H111" = 1951.4667
    "I111" = 1615
    "J111" = 2456.1667
    "K111" = 4845
    "L111" = 7368.500100000001
    "M111" = -1778
    "N111" = -13502.5001
    "H126" = 80000
    "J126" = 80000
    "L126" = 80000
    "N126" = -89880
    "H130" = 120762.5
    "J130" = 120762.5
    "L130" = 120762.5
    "N130" = -130802.5
    "H132" = 2118.5527
    "I132" = 2127.2703
    "K132" = 6381.8109
    "M132" = -3851.8109
    "H138" = 2041.27
    "I138" = 781.3333
    "J138" = 2749.9844
    "K138" = 2343.9999
    "L138" = 8249.9532
    "M138" = 2796.0001
    "N138" = -18529.9532
}
foreach ($addr in $ALC_values.Keys) {
    $ws.Range($addr).Value = $ALC_values[$addr]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ARM_values = @{
    "H32" = 43379160
    "I32" = 53053860
    "K32" = 53053860
    "M32" = -53053573
    "H61" = 2640
    "I61" = 2157.8708
    "J61" = 3636.4
    "K61" = 2157.8708
    "L61" = 3636.4
    "M61" = -1945.8708
    "N61" = -4060.4
    "H122" = 3552
    "I122" = 3281.8333
    "K122" = 9845.499899999999
    "M122" = -7395.499899999999
    "H123" = 99894.5
    "J123" = 99894.5
    "L123" = 99894.5
    "N123" = -109694.5
    "H132" = 2962.0173
    "I132" = 2311.3489
    "J132" = 4827.2666
    "K132" = 6934.0467
    "L132" = 14481.7998
    "M132" = -4404.0467
    "N132" = -19541.7998
    "H133" = 100130.5
    "J133" = 100130.5
    "L133" = 100130.5
    "N133" = -105190.5
    "H136" = 2640
    "I136" = 2157.8708
    "J136" = 3636.4
    "K136" = 6473.6124
    "L136" = 10909.2
    "M136" = -3923.6124
    "N136" = -16009.2
}
foreach ($addr in $ARM_values.Keys) {
    $ws.Range($addr).Value = $ARM_values[$addr]
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$BSM_values = @{
    "H116" = 49495
    "J116" = 49495
    "L116" = 49495
    "N116" = -58673
    "H134" = 2168384.8
    "I134" = 3108243
    "J134" = 6710.8
    "K134" = 9324729
    "L134" = 20132.4
    "M134" = -9322194
    "N134" = -25202.4
}
foreach ($addr in $BSM_values.Keys) {
    $ws.Range($addr).Value = $BSM_values[$addr]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$CRP_values = @{
    "H62" = 3764.7
    "I62" = 3390
    "K62" = 3390
    "M62" = -2766
    "H65" = 3764.7
    "I65" = 3390
    "K65" = 16950
    "M65" = -13830
    "H107" = 1200.2222
    "I107" = 576.2857
    "K107" = 576.2857
    "H122" = 3033771.8
    "I122" = 4350969.5
    "J122" = 4217.4
    "K122" = 13052908.5
    "L122" = 12652.2
    "M122" = -13050458.5
    "N122" = -17552.2
    "H132" = 3164.6765
    "I132" = 3148.0386
    "J132" = 3218.75
    "K132" = 9444.1158
    "L132" = 9656.25
    "M132" = -6914.1158
    "N132" = -14716.25
    "H134" = 4076.5
    "I134" = 4058.8572
    "K134" = 12176.5716
    "M134" = -9641.571599999999
}
foreach ($addr in $CRP_values.Keys) {
    $ws.Range($addr).Value = $CRP_values[$addr]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$CUL_values = @{
    "H140" = 2242.1
    "I140" = 1060.1428
    "K140" = 3180.4284
    "M140" = 1999.5716
}
foreach ($addr in $CUL_values.Keys) {
    $ws.Range($addr).Value = $CUL_values[$addr]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$GSM_values = @{
    "H52" = 46148.5
    "J52" = 45473
    "L52" = 45473
    "N52" = -45991
    "H113" = 10153.167
    "I113" = 2020.75
    "J113" = 26418
    "K113" = 2020.75
    "L113" = 26418
    "M113" = 149.25
    "N113" = -30758
    "H122" = 1454.4
    "I122" = 1379.25
    "J122" = 1755
    "K122" = 4137.75
    "L122" = 5265
    "M122" = -1687.75
    "N122" = -10165
    "H124" = 120999
    "J124" = 120999
    "L124" = 120999
    "N124" = -130819
    "H132" = 2579.1724
    "I132" = 2243.84
    "K132" = 6731.52
    "M132" = -4201.52
}
foreach ($addr in $GSM_values.Keys) {
    $ws.Range($addr).Value = $GSM_values[$addr]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$LTW_values = @{
    "H16" = 393.63635
    "I16" = 413
    "K16" = 413
    "M16" = -243
    "H40" = 25644928
    "I40" = 33336156
    "K40" = 33336156
    "M40" = -33336020
    "H55" = 535.7778
    "I55" = 457.83334
    "J55" = 691.6667
    "K55" = 457.83334
    "L55" = 691.6667
    "M55" = -284.83334
    "N55" = -1037.6667
    "H61" = 4683.8335
    "I61" = 0
    "J61" = 4683.8335
    "K61" = 0
    "L61" = 4683.8335
    "N61" = -5087.8335
    "H113" = 4683.8335
    "I113" = 0
    "J113" = 4683.8335
    "K113" = 0
    "L113" = 4683.8335
    "N113" = -9023.833500000001
    "H122" = 29122.072
    "I122" = 29122.072
    "K122" = 87366.216
    "M122" = -84916.216
    "H132" = 3592.8
    "I132" = 3117
    "K132" = 9351
    "M132" = -6821
    "H136" = 2888.5
    "I136" = 2684.6667
    "J136" = 3500
    "K136" = 8054.000100000001
    "L136" = 10500
    "M136" = -5504.000100000001
    "N136" = -15600
}
foreach ($addr in $LTW_values.Keys) {
    $ws.Range($addr).Value = $LTW_values[$addr]
}

$LTW_clear = @("M61", "M113")
foreach ($addr in $LTW_clear) {
    $ws.Range($addr).ClearContents()
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$WVR_values = @{
    "H32" = 11136.667
    "J32" = 10985
    "L32" = 10985
    "N32" = -11619
    "H64" = 79994.664
    "J64" = 99984
    "L64" = 99984
    "N64" = -100480
    "H67" = 79994.664
    "J67" = 99984
    "L67" = 99984
    "N67" = -101700
    "H107" = 428.05884
    "J107" = 518.25
    "L107" = 1554.75
    "N107" = -5394.75
    "H132" = 2486.4
    "I132" = 2357.75
    "K132" = 7073.25
    "M132" = -4543.25
    "H136" = 1853.8422
    "I136" = 1098.3334
    "J136" = 4687
    "K136" = 3295.0002
    "L136" = 14061
    "M136" = -745.0001999999999
    "N136" = -19161
}
foreach ($addr in $WVR_values.Keys) {
    $ws.Range($addr).Value = $WVR_values[$addr]
}

